$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin name / link), safe to set directly ---
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'

# --- Numeric-looking text cells (Price / Volume columns) ---
# Use a formula producing the literal text, then flatten the whole
# D2:E51 range to static values via copy / paste-special so Excel
# does not auto-convert the numeric-looking strings to real numbers.
$ws.Range('D2').Formula = '="51.579.20"'
$ws.Range('E2').Formula = '="  -0.63%  "'
$ws.Range('D3').Formula = '="2.792.84"'
$ws.Range('E3').Formula = '="  +0.25%  "'
$ws.Range('E4').Formula = '="  -0.01%  "'
$ws.Range('D5').Formula = '="354.00"'
$ws.Range('E5').Formula = '="  -0.76%  "'
$ws.Range('D6').Formula = '="108.90"'
$ws.Range('E6').Formula = '="  -0.59%  "'
$ws.Range('D7').Formula = '="0.555"'
$ws.Range('E7').Formula = '="  -0.85%  "'
$ws.Range('E8').Formula = '="  +0.02%  "'
$ws.Range('D9').Formula = '="0.622"'
$ws.Range('E9').Formula = '="  +5.83%  "'
$ws.Range('D10').Formula = '="39.96"'
$ws.Range('E10').Formula = '="  -0.87%  "'
$ws.Range('E11').Formula = '="  +0.95%  "'
$ws.Range('D12').Formula = '="0.0837"'
$ws.Range('E12').Formula = '="  -1.25%  "'
$ws.Range('D13').Formula = '="20.05"'
$ws.Range('E13').Formula = '="  +2.94%  "'
$ws.Range('D14').Formula = '="7.75"'
$ws.Range('E14').Formula = '="  +2.28%  "'
$ws.Range('D15').Formula = '="3.235.59"'
$ws.Range('E15').Formula = '="  +0.22%  "'
$ws.Range('D16').Formula = '="2.802.82"'
$ws.Range('E16').Formula = '="  -0.51%  "'
$ws.Range('D17').Formula = '="0.937"'
$ws.Range('E17').Formula = '="  -0.86%  "'
$ws.Range('D18').Formula = '="51.542.55"'
$ws.Range('E18').Formula = '="  -0.65%  "'
$ws.Range('D19').Formula = '="7.75"'
$ws.Range('E19').Formula = '="  +3.61%  "'
$ws.Range('D20').Formula = '="3.16"'
$ws.Range('E20').Formula = '="  +2.75%  "'
$ws.Range('D21').Formula = '="13.41"'
$ws.Range('E21').Formula = '="  +1.99%  "'
$ws.Range('D22').Formula = '="0.0₃0971"'
$ws.Range('E22').Formula = '="  -0.48%  "'
$ws.Range('D23').Formula = '="70.28"'
$ws.Range('E23').Formula = '="  +0.05%  "'
$ws.Range('D24').Formula = '="267.40"'
$ws.Range('E24').Formula = '="  -1.11%  "'
$ws.Range('D25').Formula = '="2.76"'
$ws.Range('E25').Formula = '="  +0.59%  "'
$ws.Range('D26').Formula = '="0.998"'
$ws.Range('E26').Formula = '="  -0.11%  "'
$ws.Range('D27').Formula = '="26.06"'
$ws.Range('E27').Formula = '="  -1.58%  "'
$ws.Range('E28').Formula = '="  +0.99%  "'
$ws.Range('D29').Formula = '="10.32"'
$ws.Range('E29').Formula = '="  +0.13%  "'
$ws.Range('D30').Formula = '="37.10"'
$ws.Range('E30').Formula = '="  +7.10%  "'
$ws.Range('D32').Formula = '="6.28"'
$ws.Range('E32').Formula = '="  +9.46%  "'
$ws.Range('D33').Formula = '="51.95"'
$ws.Range('E33').Formula = '="  -0.18%  "'
$ws.Range('E34').Formula = '="  +9.56%  "'
$ws.Range('D35').Formula = '="0.0443"'
$ws.Range('E35').Formula = '="  -6.00%  "'
$ws.Range('D36').Formula = '="0.0852"'
$ws.Range('E36').Formula = '="  +0.86%  "'
$ws.Range('E37').Formula = '="  -0.04%  "'
$ws.Range('D38').Formula = '="18.76"'
$ws.Range('D39').Formula = '="3.13"'
$ws.Range('E39').Formula = '="  -2.09%  "'
$ws.Range('D40').Formula = '="1.97"'
$ws.Range('E40').Formula = '="  -0.20%  "'
$ws.Range('E41').Formula = '="  +0.18%  "'
$ws.Range('D42').Formula = '="2.51"'
$ws.Range('E42').Formula = '="  -3.53%  "'
$ws.Range('D43').Formula = '="21.81"'
$ws.Range('E43').Formula = '="  -0.03%  "'
$ws.Range('D44').Formula = '="119.28"'
$ws.Range('E44').Formula = '="  -0.39%  "'
$ws.Range('D45').Formula = '="2.18"'
$ws.Range('E45').Formula = '="  -2.80%  "'
$ws.Range('D46').Formula = '="2.132.04"'
$ws.Range('E46').Formula = '="  +2.42%  "'
$ws.Range('D47').Formula = '="3.41"'
$ws.Range('E47').Formula = '="  +4.21%  "'
$ws.Range('D48').Formula = '="2.36"'
$ws.Range('E48').Formula = '="  +6.45%  "'
$ws.Range('E49').Formula = '="  +10.99%  "'
$ws.Range('D50').Formula = '="0.915"'
$ws.Range('E50').Formula = '="  -3.92%  "'
$ws.Range('D51').Formula = '="5.38"'
$ws.Range('E51').Formula = '="  -6.39%  "'

$ws.Range("D2:E51").Copy()
$ws.Range("D2:E51").PasteSpecial(-4163)
$excel.CutCopyMode = 0

